$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Activate()

# Insert a new row at 93, shifting existing rows 93..267 down to 94..268
# (this also shifts the two dataValidation ranges and every row below it).
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row with the re-added risk-group multiplier for dorms
# (it had disappeared - re-enter it right after the other riskgroup_multiplier_force_infection_* rows).
$ws.Cells.Item(93, 1).Value = "riskgroup_multiplier_force_infection_dorm"
$ws.Cells.Item(93, 2).Value = 10

# Restore the view/selection state recorded for this edit.
$ws.Range("A90").Select() | Out-Null
